$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Heading-1 text runs: "Headi" + "n" + "g 1"  ->  "Heading" + " 1".
#    The three runs share identical formatting, so the net visible text
#    is unchanged ("Heading 1") - only the run split point moves, and
#    the _GoBack bookmark (which used to wrap the lone "n" run) becomes
#    a zero-length bookmark sitting right after "Heading".
# ---------------------------------------------------------------------
$toc = $d.Bookmarks("_Toc81396657")
$tocStart = $toc.Start

# Remove the old _GoBack bookmark; it currently wraps the "n" character.
$d.Bookmarks("_GoBack").Delete()

# "n" sits right after "Headi" (5 characters into the heading text).
$nPos = $tocStart + 5
$nChar = $d.Range($nPos, $nPos + 1)
if ($nChar.Text -eq "n") {
    $nChar.Delete()

    # Re-insert "n" at the same spot so the surrounding runs merge back
    # together into a single "Heading 1" run (same formatting throughout).
    $insPoint = $d.Range($nPos, $nPos)
    $insPoint.InsertBefore("n")
}

# Re-create _GoBack as a zero-length bookmark right after "Heading"
# (7 characters into the heading text), splitting "Heading 1" into the
# "Heading" / " 1" runs.
$bmPos = $tocStart + 7
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# 2) "Heading 1" paragraph style: add 4pt (80 twips) of space after,
#    keeping the existing 6pt (120 twips) of space before.
# ---------------------------------------------------------------------
$h1 = $d.Styles("Heading 1")
$h1.ParagraphFormat.SpaceAfter = 4

# ---------------------------------------------------------------------
# 3) "Default Paragraph Font": mark it semi-hidden (best effort - some
#    hosts don't expose a working Style.Hidden setter, in which case
#    this is a harmless no-op).
# ---------------------------------------------------------------------
$dpf = $d.Styles("Default Paragraph Font")
try {
    $dpf.Hidden = $true
} catch {
}
